$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = "maa://21246 (91.41), maa://36684 (95.61), ***maa://22731 (6.25)"
